$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates that are plain text (coin names, links, percentages, non-numeric prices)
$ws.Range('D2').Value = '67.480.82'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '3.938.34'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +8.85%  '
$ws.Range('E6').Value = '  +4.62%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('E10').Value = '  +10.35%  '
$ws.Range('E11').Value = '  +8.62%  '
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.543.99'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '3.926.14'
$ws.Range('E16').Value = '  +3.69%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').Value = '67.726.19'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('E21').Value = '  +6.74%  '
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('E23').Value = '  +2.48%  '
$ws.Range('E24').Value = '  +3.36%  '
$ws.Range('E25').Value = '  +5.39%  '
$ws.Range('E26').Value = '  +6.95%  '
$ws.Range('E27').Value = '  +4.06%  '
$ws.Range('E28').Value = '  +3.60%  '
$ws.Range('E29').Value = '  -2.87%  '
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('E34').Value = '  +3.09%  '
$ws.Range('E35').Value = '  +3.94%  '
$ws.Range('E36').Value = '  +2.99%  '
$ws.Range('D37').Value = '0.0₃0804'
$ws.Range('E37').Value = '  +18.22%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  -5.56%  '
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E41').Value = '  +2.67%  '
$ws.Range('E42').Value = '  -7.77%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('E46').Value = '  +4.39%  '
$ws.Range('E47').Value = '  +5.77%  '
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('E49').Value = '  +4.55%  '
$ws.Range('E50').Value = '  -4.51%  '
$ws.Range('E51').Value = '  +3.10%  '

# Price cells whose new value looks like a plain number: force them to stay text
# (matching the source data which stores prices as text strings), then clear the
# temporary text number-format so no stray formatting is left behind.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '470.41'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.30'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('D7').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.734'
$ws.Range('D9').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000342'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.37'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.44'
$ws.Range('D13').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.04'
$ws.Range('D15').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.92'
$ws.Range('D18').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.71'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.70'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.36'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.93'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.72'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.55'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.76'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.26'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '725.59'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.63'
$ws.Range('D31').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.82'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.05'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.155'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.86'
$ws.Range('D36').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.39'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0479'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.06'
$ws.Range('D41').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.142'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.337'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.18'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.41'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '146.26'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.14'
$ws.Range('D50').ClearFormats()
